# Project "Sample Project" resave: rule row R40's "Rule" label (column B,
# row 11 on the "Rules" sheet) is changed from the text "R40" to the text
# "1". The new literal string "1" becomes a new shared-string table entry.
#
# A plain  Range.Value = "1"  assignment would be auto-coerced to a NUMBER
# by Excel's input parser (since "1" looks numeric), which would store the
# cell as a numeric literal instead of a shared string - not what we want.
# To force a genuine text cell (t="s") without leaving the cell's number
# format/style altered (target keeps the exact same style index), we:
#   1. Write a formula that evaluates to the text string "1".
#   2. Copy the cell and Paste-Special "Values only" back onto itself,
#      which collapses the formula down to its literal text result while
#      preserving the cell's existing formatting/style untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$cell = $ws.Range("B11")
$cell.Formula = "=""1"""
$cell.Copy()
$cell.PasteSpecial(-4163)   # xlPasteValues
